$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.233.88"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.858.53"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.47"
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2815"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06551"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07800"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.24"
$ws.Range("E12").Value = "  -6.93%  "
$ws.Range("D13").Value = "1.862.73"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.104"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6656"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.87"
$ws.Range("E16").Value = "  -4.23%  "
$ws.Range("D17").Value = "30.254.13"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9987"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.443"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "2.106.69"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007238"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.141"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.70"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.311"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.02"
$ws.Range("E28").Value = "  -8.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.339"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09598"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.414"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.471"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.099"
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04670"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7000"
$ws.Range("E36").Value = "  -3.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9983"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.708"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.428"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.512"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.04"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8572"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.938"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.05"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4160"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "1.004.50"
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.200"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.021"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.72"
$ws.Range("E51").Value = "  -2.84%  "
